# Add the new employee row (id=2, Thierry Hochart) at the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Thierry"
$ws.Range("C3").Value = "Hochart"

$ws.Range("C5").Select() | Out-Null
